$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common metadata columns (A-C, E-J) are identical for every row in this
# Femacal de La Calera / Chirimoya block, so set them once per row inside
# the loop below, then set the per-row varying columns (D, K, L, M, N, O,
# P, Q, R, S, T) from the data table.
# Columns per entry: 0=row, 1=Fecha, 2=Variedad, 3=Calidad, 4=Volumen,
# 5=Precio minimo, 6=Precio maximo, 7=Precio promedio ponderado,
# 8=Unidad de comercializacion, 9=Origen, 10=Precio $/Kg, 11=Kg / unidad

$rowData = @(
    @(113, 44524, "Cultivar IV Región", "Especial", 50, 26000, 26000, 26000, "$/bandeja 10 kilos", "Provincia de Limarí", 2600, 10),
    @(114, 44524, "Cultivar IV Región", "Primera", 55, 24000, 24000, 24000, "$/bandeja 10 kilos", "Provincia de Limarí", 2400, 10),
    @(115, 44524, "Cultivar IV Región", "Segunda", 50, 21000, 21000, 21000, "$/bandeja 10 kilos", "Provincia de Limarí", 2100, 10),
    @(116, 44448, "Cultivar IV Región", "Primera", 45, 30000, 30000, 30000, "$/bandeja 10 kilos", "Provincia del Elquí", 3000, 10),
    @(117, 44448, "Cultivar IV Región", "Segunda", 40, 27000, 27000, 27000, "$/bandeja 10 kilos", "Provincia del Elquí", 2700, 10),
    @(118, 44441, "Cultivar IV Región", "Primera", 68, 3000, 3000, 3000, "$/kilo (en caja de 15 kilos)", "Provincia del Elquí", 3000, 1),
    @(119, 44441, "Cultivar IV Región", "Segunda", 70, 2700, 2700, 2700, "$/kilo (en caja de 15 kilos)", "Provincia del Elquí", 2700, 1),
    @(120, 44504, "Cultivar IV Región", "Especial", 54, 26000, 26000, 26000, "$/bandeja 10 kilos", "Provincia de Limarí", 2600, 10),
    @(121, 44504, "Cultivar IV Región", "Primera", 57, 24000, 24000, 24000, "$/bandeja 10 kilos", "Provincia de Limarí", 2400, 10),
    @(122, 44504, "Cultivar IV Región", "Segunda", 50, 20000, 20000, 20000, "$/bandeja 10 kilos", "Provincia de Limarí", 2000, 10),
    @(123, 44522, "Cultivar IV Región", "Especial", 45, 26000, 26000, 26000, "$/bandeja 10 kilos", "Provincia de Limarí", 2600, 10),
    @(124, 44522, "Cultivar IV Región", "Primera", 55, 23000, 23000, 23000, "$/bandeja 10 kilos", "Provincia de Limarí", 2300, 10),
    @(125, 44522, "Cultivar IV Región", "Segunda", 45, 20000, 20000, 20000, "$/bandeja 10 kilos", "Provincia de Limarí", 2000, 10),
    @(126, 44447, "Cultivar IV Región", "Primera", 48, 30000, 30000, 30000, "$/bandeja 10 kilos", "Provincia del Elquí", 3000, 10),
    @(127, 44510, "Cultivar IV Región", "Especial", 45, 26000, 26000, 26000, "$/bandeja 10 kilos", "Provincia de Limarí", 2600, 10),
    @(128, 44510, "Cultivar IV Región", "Primera", 47, 25000, 25000, 25000, "$/bandeja 10 kilos", "Provincia de Limarí", 2500, 10),
    @(129, 44510, "Cultivar IV Región", "Segunda", 40, 21000, 21000, 21000, "$/bandeja 10 kilos", "Provincia de Limarí", 2100, 10),
    @(130, 44468, "Cultivar IV Región", "Especial", 45, 27000, 27000, 27000, "$/bandeja 10 kilos", "Provincia del Elquí", 2700, 10),
    @(131, 44468, "Cultivar IV Región", "Primera", 48, 25000, 25000, 25000, "$/bandeja 10 kilos", "Provincia del Elquí", 2500, 10),
    @(132, 44468, "Cultivar IV Región", "Segunda", 40, 22000, 22000, 22000, "$/bandeja 10 kilos", "Provincia del Elquí", 2200, 10),
    @(133, 44517, "Cultivar IV Región", "Especial", 56, 26000, 26000, 26000, "$/bandeja 10 kilos", "Provincia de Limarí", 2600, 10),
    @(134, 44517, "Cultivar IV Región", "Primera", 60, 24000, 24000, 24000, "$/bandeja 10 kilos", "Provincia de Limarí", 2400, 10),
    @(135, 44517, "Cultivar IV Región", "Segunda", 50, 22000, 22000, 22000, "$/bandeja 10 kilos", "Provincia de Limarí", 2200, 10),
    @(136, 44515, "Cultivar IV Región", "Especial", 48, 26000, 26000, 26000, "$/bandeja 10 kilos", "Provincia de Limarí", 2600, 10),
    @(137, 44515, "Cultivar IV Región", "Primera", 47, 24000, 24000, 24000, "$/bandeja 10 kilos", "Provincia de Limarí", 2400, 10),
    @(138, 44515, "Cultivar IV Región", "Segunda", 40, 22000, 22000, 22000, "$/bandeja 10 kilos", "Provincia de Limarí", 2200, 10),
    @(139, 44508, "Cultivar IV Región", "Especial", 56, 27000, 27000, 27000, "$/bandeja 10 kilos", "Provincia de Limarí", 2700, 10),
    @(140, 44508, "Cultivar IV Región", "Primera", 58, 25000, 25000, 25000, "$/bandeja 10 kilos", "Provincia de Limarí", 2500, 10),
    @(141, 44508, "Cultivar IV Región", "Segunda", 50, 22000, 22000, 22000, "$/bandeja 10 kilos", "Provincia de Limarí", 2200, 10),
    @(142, 44508, "Cultivar V Región", "Primera", 36, 24000, 24000, 24000, "$/bandeja 10 kilos", "Provincia de Quillota", 2400, 10),
    @(143, 44508, "Cultivar V Región", "Segunda", 30, 20000, 20000, 20000, "$/bandeja 10 kilos", "Provincia de Quillota", 2000, 10)
)

foreach ($entry in $rowData) {
    $r = $entry[0]

    $ws.Cells.Item($r, 1).Value = 3                            # A: Mercado ID
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"        # B: Mercado
    $ws.Cells.Item($r, 3).Value = "Coquimbo"                    # C: Region

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $entry[1]                                    # D: Fecha
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 5).Value = 5                             # E: Codreg
    $ws.Cells.Item($r, 6).Value = "Fruta"                       # F: Tipo
    $ws.Cells.Item($r, 7).Value = 100107                        # G: Producto ID
    $ws.Cells.Item($r, 8).Value = "Otros"                       # H: Producto
    $ws.Cells.Item($r, 9).Value = 100107002                     # I: Categoria ID
    $ws.Cells.Item($r, 10).Value = "Chirimoya"                  # J: Categoria

    $ws.Cells.Item($r, 11).Value = $entry[2]                    # K: Variedad
    $ws.Cells.Item($r, 12).Value = $entry[3]                    # L: Calidad
    $ws.Cells.Item($r, 13).Value = $entry[4]                    # M: Volumen
    $ws.Cells.Item($r, 14).Value = $entry[5]                    # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $entry[6]                    # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $entry[7]                    # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $entry[8]                    # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $entry[9]                    # R: Origen
    $ws.Cells.Item($r, 19).Value = $entry[10]                   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $entry[11]                   # T: Kg / unidad
}
